$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Fill in the three "总结：" (summary) placeholder cells with the full
#    summary text that was written for this iteration of the project plan.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value  = "总结：阅读完需求规格说明书后，觉得有些需求不大完善，与实际使用有些出入，所以经过小组讨论后，新增了一些实际需要的功能，例如聊天、动态分享、好友分组等，详细见用例设计文档"
$ws.Range("A19").Value = "总结：上一阶段用例设计过程出现了一些明显缺陷，例如用例名存在歧义、难以理解，功能实现过程缺少用例，出现多余用例的情况。此阶段小组成员对此进行重新讨论与修改，将有歧义的用例修正，删除多余不规范用例，添加必要用例。"
$ws.Range("A29").Value = "总结：此阶段小组编写用例规约，因为用例数量较多，所以小组成员在上一阶段修改完用例后就开始编写用例描述。注：本小组将用例描述与UI设计分成两阶段，先做用例规约，再进行UI设计"

# ---------------------------------------------------------------------------
# 2) Fill in the last (previously empty) plan table for 2018.10.11, rows
#    33-38, columns B (plan content) and C (status).
# ---------------------------------------------------------------------------
$ws.Range("B33").Value = "App端UI总体结构设计；app端地图界面设计"
$ws.Range("C33").Value = "进行中"

$ws.Range("B34").Value = "App登录界面、注册界面设计"
$ws.Range("C34").Value = "未开始"

$ws.Range("B35").Value = "App端通讯录功能界面设计"
$ws.Range("C35").Value = "进行中"

$ws.Range("B36").Value = "App端动态功能界面设计"
$ws.Range("C36").Value = "未开始"

$ws.Range("B37").Value = "App端个人功能界面设计"
$ws.Range("C37").Value = "进行中"

$ws.Range("B38").Value = "App端聊天界面设计"
$ws.Range("C38").Value = "未开始"

# The "完成情况" column (C) in that block used the plain font style while its
# neighbours (B and D) used the slightly different font flavour already
# present in the workbook; match that so the newly filled cells look the
# same as the rest of the row (copy number/font formatting only).
$fmtSource = $ws.Range("B33")
for ($r = 33; $r -le 38; $r++) {
    $fmtSource.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Reflect the scrolled / selected state the sheet was left in.
# ---------------------------------------------------------------------------
$ws.Range("A25").Select()
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B35").Select()
